$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 669.4
$ws.Range("I41").Value = 612
$ws.Range("K41").Value = 612
$ws.Range("M41").Value = -172
$ws.Range("H135").Value = 2114.923
$ws.Range("I135").Value = 1971
$ws.Range("K135").Value = 17739
$ws.Range("M135").Value = -15204
$ws.Range("H138").Value = 2159.491
$ws.Range("I138").Value = 1962.6923
$ws.Range("K138").Value = 5888.0769
$ws.Range("M138").Value = -748.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6071.9287
$ws.Range("I32").Value = 4318.364
$ws.Range("J32").Value = 12501.667
$ws.Range("K32").Value = 4318.364
$ws.Range("L32").Value = 12501.667
$ws.Range("M32").Value = -4031.364
$ws.Range("N32").Value = -13075.667
$ws.Range("H61").Value = 3366.3635
$ws.Range("I61").Value = 3374.2856
$ws.Range("K61").Value = 3374.2856
$ws.Range("M61").Value = -3162.2856
$ws.Range("H74").Value = 7473.125
$ws.Range("I74").Value = 4958
$ws.Range("J74").Value = 11665
$ws.Range("K74").Value = 4958
$ws.Range("L74").Value = 11665
$ws.Range("M74").Value = -4084
$ws.Range("N74").Value = -13413
$ws.Range("H77").Value = 7473.125
$ws.Range("I77").Value = 4958
$ws.Range("J77").Value = 11665
$ws.Range("K77").Value = 24790
$ws.Range("L77").Value = 58325
$ws.Range("M77").Value = -20422
$ws.Range("N77").Value = -67061
$ws.Range("H88").Value = 1407.0834
$ws.Range("I88").Value = 1141
$ws.Range("J88").Value = 1779.6
$ws.Range("K88").Value = 1141
$ws.Range("L88").Value = 1779.6
$ws.Range("M88").Value = -735
$ws.Range("N88").Value = -2591.6
$ws.Range("H91").Value = 1407.0834
$ws.Range("I91").Value = 1141
$ws.Range("J91").Value = 1779.6
$ws.Range("K91").Value = 1141
$ws.Range("L91").Value = 1779.6
$ws.Range("M91").Value = 263
$ws.Range("N91").Value = -4587.6
$ws.Range("H132").Value = 4015.1052
$ws.Range("I132").Value = 4587.5557
$ws.Range("J132").Value = 3499.9
$ws.Range("K132").Value = 13762.6671
$ws.Range("L132").Value = 10499.7
$ws.Range("M132").Value = -11232.6671
$ws.Range("N132").Value = -15559.7
$ws.Range("H136").Value = 3366.3635
$ws.Range("I136").Value = 3374.2856
$ws.Range("K136").Value = 10122.8568
$ws.Range("M136").Value = -7572.856800000001
$ws.Range("H140").Value = 88000
$ws.Range("J140").Value = 88000
$ws.Range("L140").Value = 88000
$ws.Range("N140").Value = -98360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3286.2
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3286.2
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 2599
$ws.Range("I99").Value = 2811.75
$ws.Range("J99").Value = 1748
$ws.Range("K99").Value = 2811.75
$ws.Range("L99").Value = 1748
$ws.Range("M99").Value = -1313.75
$ws.Range("N99").Value = -4744
$ws.Range("H134").Value = 1293.1666
$ws.Range("I134").Value = 1293.1666
$ws.Range("K134").Value = 3879.4998
$ws.Range("M134").Value = -1344.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3584.0386
$ws.Range("I31").Value = 2590.5264
$ws.Range("K31").Value = 2590.5264
$ws.Range("M31").Value = -2295.5264
$ws.Range("H34").Value = 3584.0386
$ws.Range("I34").Value = 2590.5264
$ws.Range("K34").Value = 2590.5264
$ws.Range("M34").Value = -2388.5264
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -612
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 1000
$ws.Range("K40").Value = 1000
$ws.Range("M40").Value = -840
$ws.Range("H94").Value = 4100
$ws.Range("J94").Value = 4250
$ws.Range("L94").Value = 4250
$ws.Range("N94").Value = -5152
$ws.Range("H99").Value = 14449.652
$ws.Range("I99").Value = 9959.200000000001
$ws.Range("K99").Value = 9959.200000000001
$ws.Range("M99").Value = -8461.200000000001
$ws.Range("H107").Value = 14286741
$ws.Range("I107").Value = 23810116
$ws.Range("K107").Value = 23810116
$ws.Range("M107").Value = -23808196
$ws.Range("H126").Value = 14449.652
$ws.Range("I126").Value = 9959.200000000001
$ws.Range("K126").Value = 29877.6
$ws.Range("M126").Value = -27407.6
$ws.Range("H132").Value = 2874.875
$ws.Range("I132").Value = 2499.9092
$ws.Range("J132").Value = 3699.8
$ws.Range("K132").Value = 7499.7276
$ws.Range("L132").Value = 11099.4
$ws.Range("M132").Value = -4969.7276
$ws.Range("N132").Value = -16159.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999
$ws.Range("I5").Value = 999
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2997
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2885
$ws.Range("N5").ClearContents()
$ws.Range("H34").Value = 25290
$ws.Range("I34").Value = 652.5
$ws.Range("J34").Value = 45000
$ws.Range("K34").Value = 1957.5
$ws.Range("L34").Value = 135000
$ws.Range("M34").Value = -1873.5
$ws.Range("N34").Value = -135168
$ws.Range("H38").Value = 93.666664
$ws.Range("J38").Value = 279.66666
$ws.Range("L38").Value = 838.9999799999999
$ws.Range("N38").Value = -1532.99998
$ws.Range("H51").Value = 650
$ws.Range("I51").Value = 650
$ws.Range("K51").Value = 1950
$ws.Range("M51").Value = -1490
$ws.Range("H55").Value = 28964.9
$ws.Range("I55").Value = 2216.6667
$ws.Range("J55").Value = 40428.43
$ws.Range("K55").Value = 6650.000100000001
$ws.Range("L55").Value = 121285.29
$ws.Range("M55").Value = -6473.000100000001
$ws.Range("N55").Value = -121639.29
$ws.Range("H135").Value = 999
$ws.Range("I135").Value = 999
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8991
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6456
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 7299.6665
$ws.Range("I136").Value = 949.5
$ws.Range("K136").Value = 2848.5
$ws.Range("M136").Value = 2251.5
$ws.Range("H138").Value = 1075.2858
$ws.Range("I138").Value = 1075.2858
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 3225.8574
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 1914.1426
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1399.5
$ws.Range("J126").Value = 999
$ws.Range("L126").Value = 2997
$ws.Range("N126").Value = -7937
$ws.Range("H132").Value = 1665.7778
$ws.Range("I132").Value = 1665.7778
$ws.Range("K132").Value = 4997.3334
$ws.Range("M132").Value = -2467.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 715.3333
$ws.Range("I22").Value = 715.3333
$ws.Range("K22").Value = 715.3333
$ws.Range("M22").Value = -420.3333
$ws.Range("H27").Value = 715.3333
$ws.Range("I27").Value = 715.3333
$ws.Range("K27").Value = 715.3333
$ws.Range("M27").Value = -608.3333
$ws.Range("H43").Value = 5003250
$ws.Range("I43").Value = 4003200
$ws.Range("J43").Value = 6670000
$ws.Range("K43").Value = 4003200
$ws.Range("L43").Value = 6670000
$ws.Range("M43").Value = -4003007
$ws.Range("N43").Value = -6670386
$ws.Range("H61").Value = 11114333
$ws.Range("I61").Value = 13891736
$ws.Range("K61").Value = 13891736
$ws.Range("M61").Value = -13891534
$ws.Range("H113").Value = 11114333
$ws.Range("I113").Value = 13891736
$ws.Range("K113").Value = 13891736
$ws.Range("M113").Value = -13889566
$ws.Range("H132").Value = 5497
$ws.Range("I132").Value = 4994
$ws.Range("K132").Value = 14982
$ws.Range("M132").Value = -12452
$ws.Range("H136").Value = 5203.7
$ws.Range("I136").Value = 4893.353
$ws.Range("K136").Value = 14680.059
$ws.Range("M136").Value = -12130.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 17926
$ws.Range("I11").Value = 17000
$ws.Range("J11").Value = 18234.666
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18234.666
$ws.Range("M11").Value = -16858
$ws.Range("N11").Value = -18518.666
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H107").Value = 634.2857
$ws.Range("I107").Value = 573.3333
$ws.Range("K107").Value = 1719.9999
$ws.Range("M107").Value = 200.0001
$ws.Range("H117").Value = 74999.5
$ws.Range("J117").Value = 74999.5
$ws.Range("L117").Value = 74999.5
$ws.Range("N117").Value = -84177.5
